# Commit: "Add NullFlavors for Diagnose ICD10 and function with statistics"
#
# 1) Rows 82-96 ("diagnostik_labor", "diagnostik_blutgase", ... "diagnostik_mrt"
#    block, Generation type "OPB") previously had Default values type
#    "String_test" (column D) and no Parameters/Nullflavor (columns E/F).
#    They now get:
#      D -> "String"
#      E -> "value_set=[PB, OPB]"
#      F -> "UNK"
#    Because "String_test" ends up with zero references afterwards, it drops
#    out of the shared-strings table entirely.
#
# 2) The "Diagnose ICD10 CODE" row's Parameters cell (E64) is re-entered /
#    normalised (was stored as two differently-formatted runs, now a single
#    plain run with the same text).
#
# 3) The saved view scrolls/selects to show the edited area (E64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 82; $r -le 96; $r++) {
    $ws.Range("D$r").Value = "String"
    $ws.Range("E$r").Value = "value_set=[PB, OPB]"
    $ws.Range("F$r").Value = "UNK"
}

# Re-enter the Diagnose ICD10 CODE parameters text (content unchanged, but
# this is the cell the edit session left the cursor on).
$ws.Range("E64").Value = "link=icd10gm2023.csv;column=Schlüsselnummer ohne Strich, Stern und  Ausrufezeichen"

# Leave the sheet scrolled/selected where the editor ended up.
[void]$ws.Range("A49").Select()
[void]$ws.Range("E64").Select()
